$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 19 de Mayo de 2020 a las 22:35'
$ws.Cells.Item(4, 2).Value = 1564713
$ws.Cells.Item(4, 3).Value = 14419
$ws.Cells.Item(4, 4).Value = 361465
$ws.Cells.Item(4, 5).Value = 1110114
$ws.Cells.Item(4, 7).Value = 1153
$ws.Cells.Item(4, 8).Value = 93134
$ws.Cells.Item(11, 2).Value = 177804
$ws.Cells.Item(11, 3).Value = 515
$ws.Cells.Item(11, 5).Value = 13919
$ws.Cells.Item(11, 7).Value = 62
$ws.Cells.Item(11, 8).Value = 8185
$ws.Cells.Item(55, 1).Value = 'Barein'
$ws.Cells.Item(55, 2).Value = 7532
$ws.Cells.Item(55, 3).Value = 348
$ws.Cells.Item(55, 4).Value = 2952
$ws.Cells.Item(55, 5).Value = 4568
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 12
$ws.Cells.Item(56, 1).Value = 'Argelia'
$ws.Cells.Item(56, 2).Value = 7377
$ws.Cells.Item(56, 3).Value = 176
$ws.Cells.Item(56, 4).Value = 3746
$ws.Cells.Item(56, 5).Value = 3070
$ws.Cells.Item(56, 7).Value = 6
$ws.Cells.Item(56, 8).Value = 561
$ws.Cells.Item(101, 4).Value = 91
$ws.Cells.Item(101, 5).Value = 1048
$ws.Cells.Item(134, 2).Value = 402
$ws.Cells.Item(134, 3).Value = 27
$ws.Cells.Item(134, 5).Value = 363
$ws.Cells.Item(157, 1).Value = 'Yemen'
$ws.Cells.Item(157, 2).Value = 167
$ws.Cells.Item(157, 3).Value = 37
$ws.Cells.Item(157, 4).Value = 5
$ws.Cells.Item(157, 5).Value = 134
$ws.Cells.Item(157, 7).Value = 8
$ws.Cells.Item(157, 8).Value = 28
$ws.Cells.Item(158, 1).Value = 'Guadalupe'
$ws.Cells.Item(158, 2).Value = 155
$ws.Cells.Item(158, 4).Value = 109
$ws.Cells.Item(158, 5).Value = 33
$ws.Cells.Item(158, 8).Value = 13
$ws.Cells.Item(159, 1).Value = 'Gibraltar'
$ws.Cells.Item(159, 2).Value = 147
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = 145
$ws.Cells.Item(159, 5).Value = 2
$ws.Cells.Item(160, 1).Value = 'Mozambique'
$ws.Cells.Item(160, 2).Value = 146
$ws.Cells.Item(160, 3).Value = 1
$ws.Cells.Item(160, 4).Value = 44
$ws.Cells.Item(160, 5).Value = 102
$ws.Cells.Item(160, 8).Value = 0
$ws.Cells.Item(161, 1).Value = 'Brunei'
$ws.Cells.Item(161, 2).Value = 141
$ws.Cells.Item(161, 4).Value = 136
$ws.Cells.Item(161, 5).Value = 4
$ws.Cells.Item(161, 8).Value = 1
$ws.Cells.Item(162, 1).Value = 'Mongolia'
$ws.Cells.Item(162, 2).Value = 140
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 26
$ws.Cells.Item(162, 5).Value = 114
$ws.Cells.Item(162, 8).Value = 0
$ws.Cells.Item(163, 1).Value = 'Mauritania'
$ws.Cells.Item(163, 2).Value = 131
$ws.Cells.Item(163, 3).Value = 50
$ws.Cells.Item(163, 4).Value = 7
$ws.Cells.Item(163, 5).Value = 120
$ws.Cells.Item(163, 8).Value = 4
$ws.Cells.Item(168, 1).Value = 'Islas Caimanes'
$ws.Cells.Item(168, 2).Value = 111
$ws.Cells.Item(168, 3).Value = 17
$ws.Cells.Item(168, 4).Value = 55
$ws.Cells.Item(168, 5).Value = 55
$ws.Cells.Item(168, 8).Value = 1
$ws.Cells.Item(169, 1).Value = 'Aruba'
$ws.Cells.Item(169, 2).Value = 101
$ws.Cells.Item(169, 4).Value = 93
$ws.Cells.Item(169, 5).Value = 5
$ws.Cells.Item(169, 8).Value = 3
$ws.Cells.Item(170, 1).Value = 'Monaco'
$ws.Cells.Item(170, 2).Value = 97
$ws.Cells.Item(170, 4).Value = 87
$ws.Cells.Item(170, 5).Value = 6
$ws.Cells.Item(170, 8).Value = 4
$ws.Cells.Item(171, 1).Value = 'Bahamas'
$ws.Cells.Item(171, 2).Value = 96
$ws.Cells.Item(171, 4).Value = 43
$ws.Cells.Item(171, 5).Value = 42
$ws.Cells.Item(171, 8).Value = 11
